# Export with no is_pref and no lev distance
# Rebuild the speaker-variant rows: the "id" (column B) is now always
# derived directly from the "speaker_variant" (column C) text, with no
# Levenshtein-distance grouping onto a preferred id, and the
# "is_prefered" (column D) flag is no longer populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New speaker_variant values for rows 2-15, in final row order.
$variants = @(
    "Karel",
    "Hendrik,",
    "Joris",
    "Iasper",
    "Anna Maria",
    "Mary",
    "Hendrik",
    "Kaeel",
    "Jacomo",
    "Ioris",
    "Marry,",
    "Marry",
    "Iacomo",
    "Jasper"
)

$row = 2
foreach ($variant in $variants) {
    $id = "#" + $variant.ToLower().Replace(" ", "-")

    $ws.Cells.Item($row, 2).Value = $id
    $ws.Cells.Item($row, 3).Value = $variant
    $ws.Cells.Item($row, 4).Value = $null

    $row++
}
